$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Plain value changes (no formulas) ---
$ws.Range("F12").Value = 1305498328.7699957
$ws.Range("G12").Value = 1240524717

$ws.Range("F13").Value = 325268233.58999997
$ws.Range("G13").Value = 319819483.19999999

$ws.Range("F14").Value = 44875817.619999997
$ws.Range("G14").Value = 34063116.799999997

$ws.Range("F15").Value = 35000000

$ws.Range("F16").Value = -53616441.74000001
$ws.Range("G16").Value = 60834434.380000003

$ws.Range("F19").Value = -412700000
$ws.Range("G19").Value = 3793000000

$ws.Range("F22").Value = -20015625
$ws.Range("G22").Value = 20015625

# F26 was a shared string (blank space) -> becomes a numeric value
$ws.Range("F26").Value = 1026703455.3810816
$ws.Range("G26").Value = 1029174575

# --- Formula changes: F18 and F21 become SUM formulas (previously plain values) ---
$ws.Range("F18").Formula = "=SUM(F12:F17)"
$ws.Range("F21").Formula = "=SUM(F18:F20)"

$excel.CalculateFullRebuild()
